## Annotations/Old/DavidCopperfield.xlsx -- "Add files via upload"
##
## The commit re-uploads the workbook after:
##  1. Un-hiding columns A:B on the "copperfield" sheet and giving them
##     explicit widths (they used to be hidden helper columns).
##  2. Making "copperfield" (the first sheet) the active/selected sheet
##     again instead of "Formatted" (the third sheet) - this also clears
##     the stale horizontal scroll position (topLeftCell="C1") on that
##     sheet's view.
##  3. (The absPath breadcrumb + the cached RAND() values simply follow
##     from re-saving the file from a different folder / a fresh
##     recalculation - no explicit action needed for those.)

$wb = $excel.ActiveWorkbook

# --- 1. Unhide column A and B on "copperfield" and size them ---------------
$ws = $wb.Worksheets.Item("copperfield")

$colA = $ws.Columns.Item(1)
$colB = $ws.Columns.Item(2)

$colA.Hidden = $false
$colB.Hidden = $false

$colA.ColumnWidth = 10.6
$colB.ColumnWidth = 18.5

# --- 2. Select "copperfield" as the active sheet again ----------------------
$ws.Activate()
$ws.Range("I267").Select()
